$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextCell 'D2' '37.332.08'
Set-TextCell 'E2' '  +2.55%  '

Set-TextCell 'D3' '2.064.59'
Set-TextCell 'E3' '  +3.93%  '

Set-TextCell 'E4' '  +0.03%  '

Set-TextCell 'D5' '234.89'
Set-TextCell 'E5' '  -0.20%  '

Set-TextCell 'E6' '  +2.85%  '

Set-TextCell 'D7' '58.09'
Set-TextCell 'E7' '  +6.90%  '

Set-TextCell 'E8' '  +0.04%  '

Set-TextCell 'D9' '0.383'
Set-TextCell 'E9' '  +3.68%  '

Set-TextCell 'D10' '59.03'
Set-TextCell 'E10' '  +2.00%  '

Set-TextCell 'D11' '0.0761'
Set-TextCell 'E11' '  +2.00%  '

Set-TextCell 'E12' '  +3.15%  '

Set-TextCell 'D13' '2.369.71'

Set-TextCell 'D14' '14.59'
Set-TextCell 'E14' '  +3.05%  '

Set-TextCell 'D15' '20.95'
Set-TextCell 'E15' '  +3.72%  '

Set-TextCell 'D16' '0.778'
Set-TextCell 'E16' '  +3.22%  '

Set-TextCell 'D17' '5.19'
Set-TextCell 'E17' '  +3.14%  '

Set-TextCell 'D18' '2.073.18'
Set-TextCell 'E18' '  +4.52%  '

Set-TextCell 'D19' '37.596.55'
Set-TextCell 'E19' '  +3.32%  '

Set-TextCell 'D20' '6.17'
Set-TextCell 'E20' '  +17.35%  '

Set-TextCell 'D21' '69.00'
Set-TextCell 'E21' '  +1.95%  '

Set-TextCell 'D22' '0.0₃0814'
Set-TextCell 'E22' '  +1.58%  '

Set-TextCell 'D23' '226.26'
Set-TextCell 'E23' '  +2.32%  '

Set-TextCell 'D24' '0.999'
Set-TextCell 'E24' '  -0.02%  '

Set-TextCell 'D25' '2.45'
Set-TextCell 'E25' '  +2.18%  '

Set-TextCell 'D26' '2.38'
Set-TextCell 'E26' '  +1.23%  '

Set-TextCell 'D27' '164.49'
Set-TextCell 'E27' '  +1.13%  '

Set-TextCell 'E28' '  +13.68%  '

Set-TextCell 'D29' '8.87'
Set-TextCell 'E29' '  +2.38%  '

Set-TextCell 'D30' '19.17'
Set-TextCell 'E30' '  +1.84%  '

Set-TextCell 'D31' '0.126'
Set-TextCell 'E31' '  -1.25%  '

Set-TextCell 'E32' '  +1.81%  '

Set-TextCell 'D33' '4.49'
Set-TextCell 'E33' '  +2.75%  '

Set-TextCell 'D34' '0.0621'
Set-TextCell 'E34' '  +2.88%  '

Set-TextCell 'E35' '  +8.90%  '

Set-TextCell 'E36' '  +6.22%  '

Set-TextCell 'D37' '3.39'
Set-TextCell 'E37' '  +1.94%  '

Set-TextCell 'E38' '  +0.18%  '

Set-TextCell 'E39' '  +0.72%  '

Set-TextCell 'D40' '5.86'
Set-TextCell 'E40' '  +7.03%  '

Set-TextCell 'E41' '  +7.06%  '

Set-TextCell 'E42' '  -1.34%  '

Set-TextCell 'D43' '1.474.23'
Set-TextCell 'E43' '  +1.47%  '

Set-TextCell 'D44' '96.77'
Set-TextCell 'E44' '  +8.29%  '

Set-TextCell 'B45' 'TrustWalletToken'
Set-TextCell 'C45' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 'D45' '1.16'
Set-TextCell 'E45' '  +5.94%  '

Set-TextCell 'B46' 'FTXToken'
Set-TextCell 'C46' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell 'D46' '4.28'
Set-TextCell 'E46' '  +18.89%  '

Set-TextCell 'E47' '  +4.25%  '

Set-TextCell 'D48' '15.94'
Set-TextCell 'E48' '  +6.65%  '

Set-TextCell 'E49' '  +4.02%  '

Set-TextCell 'D50' '7.27'
Set-TextCell 'E50' '  +6.61%  '

Set-TextCell 'E51' '  +2.25%  '
